$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new header cells (I1 = "I0", J1 = "IF"), copying the
# existing bold/bordered/centered header style from H1 so the new
# header cells match the look of the other headers.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Row-by-row data for the new I (I0) and J (IF) columns.
$data = @(
    @(2, 1, 6),
    @(3, 1, 5),
    @(4, 1, 5),
    @(5, 1, 5),
    @(6, 1, 7),
    @(7, 1, 4),
    @(8, 1, 5),
    @(9, 1, 5),
    @(10, 2, 6),
    @(11, 1, 4),
    @(12, 8, 8),
    @(13, 7, 8),
    @(14, 6, 6),
    @(15, 4, 5),
    @(16, 4, 5),
    @(17, 6, 7),
    @(18, 6, 7),
    @(19, 6, 6),
    @(20, 7, 8),
    @(21, 6, 7),
    @(22, 6, 6),
    @(23, 6, 7),
    @(24, 8, 8),
    @(25, 5, 6),
    @(26, 4, 6),
    @(27, 5, 6),
    @(28, 8, 8),
    @(29, 7, 7),
    @(30, 4, 5),
    @(31, 8, 8),
    @(32, 1, 1),
    @(33, 6, 6),
    @(34, 6, 6),
    @(35, 7, 7),
    @(36, 8, 8),
    @(37, 7, 7)
)

foreach ($row in $data) {
    $r = $row[0]
    $iVal = $row[1]
    $jVal = $row[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}

Write-Output "Added I0 and IF columns"
